# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# list with freshly scraped values, written as plain text (matching the
# original sheet which stores these figures as strings, not numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.945.11'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +3.09%  '
$ws.Range("E2").Style = "Normal"
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.049.46'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +5.57%  '
$ws.Range("E3").Style = "Normal"
# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.36%  '
$ws.Range("E4").Style = "Normal"
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '512.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +4.90%  '
$ws.Range("E5").Style = "Normal"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.19'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +5.67%  '
$ws.Range("E6").Style = "Normal"
# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E7").Style = "Normal"
# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.83%  '
$ws.Range("E8").Style = "Normal"
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.17'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.00%  '
$ws.Range("E9").Style = "Normal"
# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +4.37%  '
$ws.Range("E10").Style = "Normal"
# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +6.10%  '
$ws.Range("E11").Style = "Normal"
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.574.12'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +5.97%  '
$ws.Range("E12").Style = "Normal"
# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.07%  '
$ws.Range("E13").Style = "Normal"
# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.11%  '
$ws.Range("E14").Style = "Normal"
# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.98%  '
$ws.Range("E15").Style = "Normal"
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '57.014.55'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.39%  '
$ws.Range("E16").Style = "Normal"
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.051.50'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +5.95%  '
$ws.Range("E17").Style = "Normal"
# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.72%  '
$ws.Range("E18").Style = "Normal"
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.01'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +5.60%  '
$ws.Range("E19").Style = "Normal"
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.09'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +6.73%  '
$ws.Range("E20").Style = "Normal"
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '334.26'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +7.27%  '
$ws.Range("E21").Style = "Normal"
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.31%  '
$ws.Range("E22").Style = "Normal"
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.501'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +5.18%  '
$ws.Range("E23").Style = "Normal"
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.07'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +5.14%  '
$ws.Range("E24").Style = "Normal"
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.167'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +5.52%  '
$ws.Range("E25").Style = "Normal"
# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.20%  '
$ws.Range("E26").Style = "Normal"
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0₃0930'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +11.42%  '
$ws.Range("E27").Style = "Normal"
# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.40%  '
$ws.Range("E28").Style = "Normal"
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.88'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.82%  '
$ws.Range("E29").Style = "Normal"
# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.33%  '
$ws.Range("E30").Style = "Normal"
# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +5.33%  '
$ws.Range("E31").Style = "Normal"
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.17'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +5.05%  '
$ws.Range("E32").Style = "Normal"
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '154.55'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.20%  '
$ws.Range("E33").Style = "Normal"
# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.24%  '
$ws.Range("E34").Style = "Normal"
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.81'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +5.21%  '
$ws.Range("E35").Style = "Normal"
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '26.33'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +8.13%  '
$ws.Range("E36").Style = "Normal"
# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +4.00%  '
$ws.Range("E37").Style = "Normal"
# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.03%  '
$ws.Range("E38").Style = "Normal"
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.085.07'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +5.94%  '
$ws.Range("E39").Style = "Normal"
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.83'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.19%  '
$ws.Range("E40").Style = "Normal"
# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.51%  '
$ws.Range("E41").Style = "Normal"
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.666'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +6.63%  '
$ws.Range("E42").Style = "Normal"
# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +5.47%  '
$ws.Range("E43").Style = "Normal"
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.231.78'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +7.33%  '
$ws.Range("E44").Style = "Normal"
# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +9.28%  '
$ws.Range("E45").Style = "Normal"
# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.97%  '
$ws.Range("E46").Style = "Normal"
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.931'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.46%  '
$ws.Range("E47").Style = "Normal"
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.81'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.20%  '
$ws.Range("E48").Style = "Normal"
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.69'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +6.35%  '
$ws.Range("E49").Style = "Normal"
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0865'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.63%  '
$ws.Range("E50").Style = "Normal"
# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +5.51%  '
$ws.Range("E51").Style = "Normal"

Write-Output "Applied cryptos list update"
